$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4 (shifts Calça/Vestido down) and insert "Jaqueta"
$ws.Rows.Item(4).Insert()
$ws.Range("A4").Value = "Jaqueta"

# Append the remaining new products after the existing list (now ending at row 6)
$ws.Range("A7").Value = "Bermuda"
$ws.Range("A8").Value = "Tênis"
$ws.Range("A9").Value = "Bolsa"
$ws.Range("A10").Value = "Boné"
$ws.Range("A11").Value = "Cinto"

# Add the new "Tamanho" header in column B
$ws.Range("B1").Value = "Tamanho"

$ws.Range("A11").Select()
